# Repull data, push all data, mean calculation
# Update the "dSF" column (F) values with repulled data for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = 7
$ws.Range("F5").Value = -7
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 0
$ws.Range("F14").Value = -9
$ws.Range("F15").Value = 6
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = -1
